# Rename the "SwateTemplateMetadata" sheet to "isa_template"
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("SwateTemplateMetadata")
$ws.Name = "isa_template"
